# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.807.73"
$ws.Range("E2").Value = "'  +1.80%  "

$ws.Range("D3").Value = "'2.029.36"
$ws.Range("E3").Value = "'  +0.94%  "

$ws.Range("E4").Value = "'  +0.06%  "

$ws.Range("D5").Value = "'249.33"
$ws.Range("E5").Value = "'  -0.91%  "

$ws.Range("E6").Value = "'  -0.90%  "

$ws.Range("D7").Value = "'63.09"
$ws.Range("E7").Value = "'  +0.96%  "

$ws.Range("E8").Value = "'  +0.03%  "

$ws.Range("D9").Value = "'0.393"
$ws.Range("E9").Value = "'  +6.88%  "

$ws.Range("D10").Value = "'58.25"
$ws.Range("E10").Value = "'  -1.35%  "

$ws.Range("D11").Value = "'0.0800"
$ws.Range("E11").Value = "'  +7.36%  "

$ws.Range("D13").Value = "'0.900"
$ws.Range("E13").Value = "'  -0.51%  "

$ws.Range("D14").Value = "'23.55"
$ws.Range("E14").Value = "'  +20.01%  "

$ws.Range("E15").Value = "'  -3.01%  "

$ws.Range("D16").Value = "'2.327.12"
$ws.Range("E16").Value = "'  +0.93%  "

$ws.Range("D17").Value = "'5.57"
$ws.Range("E17").Value = "'  +2.58%  "

$ws.Range("D18").Value = "'2.025.07"
$ws.Range("E18").Value = "'  +1.04%  "

$ws.Range("D19").Value = "'36.788.56"
$ws.Range("E19").Value = "'  +1.95%  "

$ws.Range("D20").Value = "'72.45"
$ws.Range("E20").Value = "'  +0.62%  "

$ws.Range("D21").Value = "'0.0₃0888"
$ws.Range("E21").Value = "'  +3.79%  "

$ws.Range("D22").Value = "'5.40"
$ws.Range("E22").Value = "'  +2.77%  "

$ws.Range("D23").Value = "'236.86"
$ws.Range("E23").Value = "'  +1.39%  "

$ws.Range("E24").Value = "'  +0.09%  "

$ws.Range("D25").Value = "'2.54"
$ws.Range("E25").Value = "'  -4.50%  "

$ws.Range("E26").Value = "'  +1.00%  "

$ws.Range("D27").Value = "'9.92"
$ws.Range("E27").Value = "'  +4.38%  "

$ws.Range("E28").Value = "'  +26.26%  "

$ws.Range("B29").Value = "'Monero"
$ws.Range("C29").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'160.51"
$ws.Range("E29").Value = "'  -2.11%  "

$ws.Range("B30").Value = "'EthereumClassic"
$ws.Range("C30").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'20.45"
$ws.Range("E30").Value = "'  +4.58%  "

$ws.Range("E31").Value = "'  +0.72%  "

$ws.Range("E32").Value = "'  -0.51%  "

$ws.Range("D33").Value = "'1.19"
$ws.Range("E33").Value = "'  +0.21%  "

$ws.Range("E34").Value = "'  +3.36%  "

$ws.Range("D35").Value = "'4.50"
$ws.Range("E35").Value = "'  -0.20%  "

$ws.Range("B36").Value = "'THORChain"
$ws.Range("C36").Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D36").Value = "'6.51"
$ws.Range("E36").Value = "'  +11.70%  "

$ws.Range("B37").Value = "'LidoDAOToken"
$ws.Range("C37").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'2.40"
$ws.Range("E37").Value = "'  -4.08%  "

$ws.Range("E38").Value = "'  +0.08%  "

$ws.Range("E39").Value = "'  +1.04%  "

$ws.Range("E40").Value = "'  +31.84%  "

$ws.Range("B41").Value = "'Cronos"
$ws.Range("C41").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").Value = "'0.100"
$ws.Range("E41").Value = "'  -3.26%  "

$ws.Range("B42").Value = "'TrustWalletToken"
$ws.Range("C42").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.26"
$ws.Range("E42").Value = "'  +3.21%  "

$ws.Range("E43").Value = "'  +1.67%  "

$ws.Range("D45").Value = "'17.03"
$ws.Range("E45").Value = "'  +2.49%  "

$ws.Range("E46").Value = "'  +0.23%  "

$ws.Range("D47").Value = "'93.92"
$ws.Range("E47").Value = "'  +0.20%  "

$ws.Range("D48").Value = "'7.67"
$ws.Range("E48").Value = "'  -1.63%  "

$ws.Range("D49").Value = "'1.373.01"
$ws.Range("E49").Value = "'  -3.38%  "

$ws.Range("D50").Value = "'2.91"
$ws.Range("E50").Value = "'  -0.02%  "

$ws.Range("D51").Value = "'2.214.50"
$ws.Range("E51").Value = "'  +0.94%  "
